# Update marksheet correct/total marks.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row: marks awarded per correct answer (was 3, now 5)
$ws.Range("B11").Value = 5

# "Total" row: total marks scored (Right * Marking) and Score/Max string
$ws.Range("B12").Value = 95
$ws.Range("E12").Value = "95/140"
